$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values per column (B..Q), applied uniformly to rows 2-26
$values = @{
    "B" = 0.5692387035740197
    "C" = 0.2355707868591003
    "D" = 0.809700823863239
    "E" = -1.520158467522549
    "F" = 0.5640702052195777
    "G" = 0.2557185134775685
    "H" = 0.4537982025430372
    "I" = 0.2157781251675712
    "J" = 0.3391897445557057
    "K" = 0.2774839348616385
    "L" = 0.2823828768330753
    "M" = 0.5056861808251918
    "N" = 0.06015717143422494
    "O" = 0.5272142884576241
    "P" = 28.72735599237364
    "Q" = 44.57274171566024
}

foreach ($row in 2..26) {
    foreach ($col in $values.Keys) {
        $ws.Range("$col$row").Value = $values[$col]
    }
}
